$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 916.86206
$ws.Range("J17").Value = 694.7143
$ws.Range("L17").Value = 2084.1429
$ws.Range("N17").Value = -2420.1429
$ws.Range("H43").Value = 2409.125
$ws.Range("J43").Value = 2815.3333
$ws.Range("L43").Value = 2815.3333
$ws.Range("N43").Value = -2953.3333
$ws.Range("H75").Value = 34666.668
$ws.Range("J75").Value = 34666.668
$ws.Range("L75").Value = 34666.668
$ws.Range("N75").Value = -36538.668
$ws.Range("H78").Value = 34666.668
$ws.Range("J78").Value = 34666.668
$ws.Range("L78").Value = 104000.004
$ws.Range("N78").Value = -113360.004
$ws.Range("H93").Value = 32675.812
$ws.Range("J93").Value = 32675.812
$ws.Range("L93").Value = 32675.812
$ws.Range("N93").Value = -37667.81200000001
$ws.Range("H123").Value = 41890
$ws.Range("J123").Value = 41890
$ws.Range("L123").Value = 41890
$ws.Range("N123").Value = -51690
$ws.Range("H125").Value = 1127.1538
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 1113.909
$ws.Range("K125").Value = 10800
$ws.Range("L125").Value = 10025.181
$ws.Range("M125").Value = -8340
$ws.Range("N125").Value = -14945.181
$ws.Range("H137").Value = 2937.647
$ws.Range("I137").Value = 2205.4546
$ws.Range("J137").Value = 4280
$ws.Range("K137").Value = 6616.3638
$ws.Range("L137").Value = 12840
$ws.Range("M137").Value = -4066.3638
$ws.Range("N137").Value = -17940
$ws.Range("H138").Value = 3595.82
$ws.Range("I138").Value = 702.8461
$ws.Range("J138").Value = 4612.2705
$ws.Range("K138").Value = 2108.5383
$ws.Range("L138").Value = 13836.8115
$ws.Range("M138").Value = 3031.4617
$ws.Range("N138").Value = -24116.8115
$ws.Range("H139").Value = 33951.613
$ws.Range("J139").Value = 33951.613
$ws.Range("L139").Value = 33951.613
$ws.Range("N139").Value = -44231.613
$ws.Range("H141").Value = 47330.5
$ws.Range("I141").Value = 51346.91
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 154040.73
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = -148860.73
$ws.Range("N141").Value = -19810

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4818.574
$ws.Range("I32").Value = 4328.477
$ws.Range("J32").Value = 6975
$ws.Range("K32").Value = 4328.477
$ws.Range("L32").Value = 6975
$ws.Range("M32").Value = -4041.477
$ws.Range("N32").Value = -7549
$ws.Range("H61").Value = 1321.8182
$ws.Range("I61").Value = 1423.2727
$ws.Range("J61").Value = 1118.909
$ws.Range("K61").Value = 1423.2727
$ws.Range("L61").Value = 1118.909
$ws.Range("M61").Value = -1211.2727
$ws.Range("N61").Value = -1542.909
$ws.Range("H136").Value = 1321.8182
$ws.Range("I136").Value = 1423.2727
$ws.Range("J136").Value = 1118.909
$ws.Range("K136").Value = 4269.8181
$ws.Range("L136").Value = 3356.727
$ws.Range("M136").Value = -1719.8181
$ws.Range("N136").Value = -8456.727000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2248.3125
$ws.Range("I86").Value = 2067
$ws.Range("K86").Value = 2067
$ws.Range("M86").Value = -944
$ws.Range("H89").Value = 2248.3125
$ws.Range("I89").Value = 2067
$ws.Range("K89").Value = 10335
$ws.Range("M89").Value = -4719

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1916.4407
$ws.Range("I58").Value = 1747.8868
$ws.Range("J58").Value = 3405.3333
$ws.Range("K58").Value = 1747.8868
$ws.Range("L58").Value = 3405.3333
$ws.Range("M58").Value = -1544.8868
$ws.Range("N58").Value = -3811.3333
$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812
$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808
$ws.Range("H134").Value = 29259
$ws.Range("I134").Value = 34012
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 102036
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -99501
$ws.Range("N134").Value = -50070
$ws.Range("H136").Value = 1916.4407
$ws.Range("I136").Value = 1747.8868
$ws.Range("J136").Value = 3405.3333
$ws.Range("K136").Value = 5243.6604
$ws.Range("L136").Value = 10215.9999
$ws.Range("M136").Value = -2693.6604
$ws.Range("N136").Value = -15315.9999
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 364.85715
$ws.Range("I92").Value = 343.75
$ws.Range("J92").Value = 393
$ws.Range("K92").Value = 1031.25
$ws.Range("L92").Value = 1179
$ws.Range("M92").Value = 216.75
$ws.Range("N92").Value = -3675
$ws.Range("H113").Value = 679.4888999999999
$ws.Range("I113").Value = 680.25
$ws.Range("J113").Value = 677.61536
$ws.Range("K113").Value = 2040.75
$ws.Range("L113").Value = 2032.84608
$ws.Range("M113").Value = 129.25
$ws.Range("N113").Value = -6372.84608
$ws.Range("H129").Value = 3292.5833
$ws.Range("J129").Value = 4171
$ws.Range("L129").Value = 12513
$ws.Range("N129").Value = -22513
$ws.Range("H131").Value = 5814777.5
$ws.Range("J131").Value = 845.7590300000001
$ws.Range("L131").Value = 2537.27709
$ws.Range("N131").Value = -12617.27709

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 19999
$ws.Range("J3").Value = 19999
$ws.Range("L3").Value = 19999
$ws.Range("N3").Value = -20223
$ws.Range("H15").Value = 19999
$ws.Range("J15").Value = 19999
$ws.Range("L15").Value = 19999
$ws.Range("N15").Value = -20339
$ws.Range("H20").Value = 15500
$ws.Range("J20").Value = 15500
$ws.Range("L20").Value = 15500
$ws.Range("N20").Value = -15952
$ws.Range("H22").Value = 85592.414
$ws.Range("I22").Value = 201540.2
$ws.Range("J22").Value = 2772.5715
$ws.Range("K22").Value = 201540.2
$ws.Range("L22").Value = 2772.5715
$ws.Range("M22").Value = -201245.2
$ws.Range("N22").Value = -3362.5715
$ws.Range("H27").Value = 85592.414
$ws.Range("I27").Value = 201540.2
$ws.Range("J27").Value = 2772.5715
$ws.Range("K27").Value = 201540.2
$ws.Range("L27").Value = 2772.5715
$ws.Range("M27").Value = -201433.2
$ws.Range("N27").Value = -2986.5715
$ws.Range("H46").Value = 3337.5
$ws.Range("I46").Value = 3666.6667
$ws.Range("J46").Value = 3140
$ws.Range("K46").Value = 3666.6667
$ws.Range("L46").Value = 3140
$ws.Range("M46").Value = -3478.6667
$ws.Range("N46").Value = -3516
$ws.Range("H55").Value = 649.6667
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 690.85
$ws.Range("I68").Value = 690.85
$ws.Range("K68").Value = 690.85
$ws.Range("M68").Value = 58.14999999999998
$ws.Range("H71").Value = 690.85
$ws.Range("I71").Value = 690.85
$ws.Range("K71").Value = 3454.25
$ws.Range("M71").Value = 289.75
$ws.Range("H132").Value = 5009.48
$ws.Range("I132").Value = 1747.4286
$ws.Range("J132").Value = 9161.182000000001
$ws.Range("K132").Value = 5242.2858
$ws.Range("L132").Value = 27483.546
$ws.Range("M132").Value = -2712.2858
$ws.Range("N132").Value = -32543.546
$ws.Range("H133").Value = 33720
$ws.Range("J133").Value = 33720
$ws.Range("L133").Value = 33720
$ws.Range("N133").Value = -38780

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 725.2857
$ws.Range("I107").Value = 762.8333
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 2288.4999
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -368.4998999999998
$ws.Range("N107").Value = -5340
$ws.Range("H136").Value = 12867.833
$ws.Range("I136").Value = 9105.5
$ws.Range("J136").Value = 14749
$ws.Range("K136").Value = 27316.5
$ws.Range("L136").Value = 44247
$ws.Range("M136").Value = -24766.5
$ws.Range("N136").Value = -49347
$ws.Range("H139").Value = 37965.383
$ws.Range("I139").Value = 40650
$ws.Range("J139").Value = 37741.668
$ws.Range("K139").Value = 40650
$ws.Range("L139").Value = 37741.668
$ws.Range("M139").Value = -35510
$ws.Range("N139").Value = -48021.668
$ws.Range("H141").Value = 42795.91
$ws.Range("J141").Value = 42795.91
$ws.Range("L141").Value = 42795.91
$ws.Range("N141").Value = -53155.91
